$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "error_Messages_Returned" column (old column H) entirely;
# this shifts the old I/J/K columns (comments, X_COORD, Y_COORD) left to H/I/J.
$ws.Columns("H").Delete()

# The X/Y coordinate pair on row 2 becomes a different (text) value pair.
$ws.Range("J2").Value = "1000"
$ws.Range("I2").Value = "800"

# Update the trade-location-verification text for rows 3 and 4 to reflect
# different option selections ("Give Feedback" iteration).
$ws.Range("G3").Value = "Crown Thinning in PO12 2HE"
$ws.Range("G4").Value = "Electric Boiler Repairs in Orkney"

# Widen the (now-relocated) tradeLocationVerification column to fit the longer text.
$ws.Columns("G").ColumnWidth = 47.140625

# Match the new selection left behind on the sheet (columns H:I, full height).
$ws.Range("H1:I1048576").Select()
